$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.580.57"
$ws.Range("E2").Value = "  +2.40%  "

$ws.Range("D3").Value = "'1.875.58"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("E4").Value = "  +0.63%  "

$ws.Range("D5").Value = "'315.96"
$ws.Range("E5").Value = "  +1.82%  "

$ws.Range("E6").Value = "  +1.09%  "

$ws.Range("D7").Value = "'0.5102"
$ws.Range("E7").Value = "  +1.05%  "

$ws.Range("D8").Value = "'0.3927"
$ws.Range("E8").Value = "  +1.77%  "

$ws.Range("D9").Value = "'0.08407"
$ws.Range("E9").Value = "  +3.46%  "

$ws.Range("E10").Value = "  +0.94%  "

$ws.Range("D11").Value = "'41.78"
$ws.Range("E11").Value = "  +1.33%  "

$ws.Range("D12").Value = "'6.272"
$ws.Range("E12").Value = "  +2.42%  "

$ws.Range("D13").Value = "'1.879.98"
$ws.Range("E13").Value = "  +2.88%  "

$ws.Range("E14").Value = "  +2.23%  "

$ws.Range("D15").Value = "'7.266"
$ws.Range("E15").Value = "  +1.89%  "

$ws.Range("D16").Value = "'1.008"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").Value = "'91.44"
$ws.Range("E18").Value = "  +1.54%  "

$ws.Range("D19").Value = "'0.06735"
$ws.Range("E19").Value = "  +1.76%  "

$ws.Range("E20").Value = "  +1.45%  "

$ws.Range("E21").Value = "  +1.23%  "

$ws.Range("D22").Value = "'5.965"
$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").Value = "'28.607.02"
$ws.Range("E23").Value = "  +2.46%  "

$ws.Range("D24").Value = "'11.14"
$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("D25").Value = "'2.248"
$ws.Range("E25").Value = "  +1.81%  "

$ws.Range("D26").Value = "'2.090.09"
$ws.Range("E26").Value = "  +2.49%  "

$ws.Range("D27").Value = "'161.86"
$ws.Range("E27").Value = "  +2.36%  "

$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").Value = "'2.366"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").Value = "'126.67"
$ws.Range("E30").Value = "  +1.81%  "

$ws.Range("D31").Value = "'0.1056"
$ws.Range("E31").Value = "  +1.91%  "

$ws.Range("E32").Value = "  +2.70%  "

$ws.Range("D33").Value = "'5.805"
$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("D34").Value = "'3.612"
$ws.Range("E34").Value = "  +0.85%  "

$ws.Range("E35").Value = "  +2.59%  "

$ws.Range("D36").Value = "'0.06526"
$ws.Range("E36").Value = "  +0.96%  "

$ws.Range("E37").Value = "  +1.21%  "

$ws.Range("E38").Value = "  -3.09%  "

$ws.Range("D39").Value = "'1.266"
$ws.Range("E39").Value = "  +3.59%  "

$ws.Range("E40").Value = "  +2.93%  "

$ws.Range("D41").Value = "'0.6465"
$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("E42").Value = "  +3.07%  "

$ws.Range("E43").Value = "  +1.38%  "

$ws.Range("E44").Value = "  +1.19%  "

$ws.Range("D45").Value = "'0.6067"
$ws.Range("E45").Value = "  +1.27%  "

$ws.Range("D46").Value = "'12.99"
$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("D47").Value = "'3.702"
$ws.Range("E47").Value = "  +1.85%  "

$ws.Range("E48").Value = "  +3.14%  "

$ws.Range("D49").Value = "'1.218"
$ws.Range("E49").Value = "  +2.38%  "

$ws.Range("D50").Value = "'122.46"
$ws.Range("E50").Value = "  +1.80%  "

$ws.Range("D51").Value = "'1.197"
$ws.Range("E51").Value = "  -5.40%  "
